$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New small side table (H4:I5) moved from old A11/A12/A13 area ---
$ws.Range("H4").Value = "QBOARD BASIQ Bouwplaat 1200mm x 600mm"
$ws.Range("I4").Value = 13.45

$ws.Range("H5").Value = "QBOARD boardfix montage lijm 290 ml"
$ws.Range("I5").Value = 12.91

# --- Fill out rows 10-13 in columns A/B with the new expense items ---
$ws.Range("A10").Value = "Opbergbox voor demonstratie"
$ws.Range("B10").Value = 3.82

$ws.Range("A11").Value = "Afdekfolie 4x5m"
$ws.Range("B11").Value = 0.65

$ws.Range("A12").Value = "Houten plaat, prototype 2"
$ws.Range("B12").Value = 5

$ws.Range("A13").Value = "Arcylplaat 25x50cm"
$ws.Range("B13").Value = 12

# --- Row 2: new "concept uitgaves" label in H2 ---
$ws.Range("H2").Value = "concept uitgaves"

# --- Set width for the newly used column H ---
$ws.Columns.Item(8).ColumnWidth = 36.67

# --- Update active selection to match the new sheet state ---
$ws.Range("B15").Select() | Out-Null
